# Ajustes realizados durante la documentación
#
# La fila de "Marta" (fila 6) tenía varios datos erróneos que se corrigen
# aquí: apellidos, DNI y edad. El email se actualiza a una nueva dirección
# y se convierte en un hipervínculo "mailto:".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apellidos: "Pérez Gómez" -> "Alcaide Perez"
$ws.Range("B6").Value = "Alcaide Perez"

# DNI: "33445566E" -> "33445566Z"
$ws.Range("C6").Value = "33445566Z"

# Edad: 29 -> 56
$ws.Range("D6").Value = 56

# Email: "marta@iem.com" -> "martanuevo@iem.com", convertido en hipervínculo
$ws.Range("E6").Value = "martanuevo@iem.com"
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:martanuevo@iem.com") | Out-Null

# Deja la selección tal y como quedó guardada en el libro final
$ws.Range("G6").Select() | Out-Null
